# Quarterly database update: roll the reporting window forward by one
# quarter (drop 1399/06, add 1401/12) and refresh historical figures per
# the updated read_price algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header quarter labels (row 8 and row 24), columns E..N ---------------
$quarters = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

for ($i = 0; $i -lt $quarters.Length; $i++) {
    $col = 5 + $i   # E=5 .. N=14
    $ws.Cells.Item(8, $col).Value = $quarters[$i]
    $ws.Cells.Item(24, $col).Value = $quarters[$i]
}

# --- Data rows, columns E..N ------------------------------------------------
$rows = @{
    10 = @(1060038, 943970, 1446441, 870019, 1532072, 807816, 925241, 363622, 3124859, 3443875)
    11 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    12 = @(100620, 205313, 42243, 677760, 339690, 335132, 299947, 203331, 371112, 393451)
    13 = @(494858, 413495, 32142, 1115897, 970830, 1028127, 404190, 112584, 3248226, 1305000)
    14 = @(-696714, 0, 0, 0, 0, 0, 0, 1813451, -1813451, 0)
    15 = @(-84026, 0, 0, 0, 0, 0, 0, 488624, -488624, 0)
    16 = @(132391, 56543, 65760, 100776, 59908, 99925, 106292, 148543, 243694, 290754)
    17 = @(1469192, 2184916, 1753756, 3871887, 2417332, 2874754, 4315720, 4232462, 3249600, 5098284)
    19 = @(2585817, 2552404, 2976785, 1433666, 2421157, 7053866, 5634772, 1852034, 6864894, 7074800)
    20 = @(5062176, 6356641, 6317127, 8070005, 7740989, 12199620, 11686162, 9214651, 14800310, 17606164)
    26 = @(2269, 2230, 1764, 1752, 2573, 2534, 2534, 3326, 3417, 1679)
    27 = @(10214, 9999, 10441, 10581, 9721, 9627, 9627, 8787, 8698, 10267)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i   # E=5 .. N=14
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
